$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new cells for the administrator row (row 2)
$ws.Range("B2").Value = "管理员的借阅记录分析"
$ws.Range("C2").Value = "管理员的系统维护"

# Rename the reader's analysis/maintenance entries in row 3 (D3/E3)
$ws.Range("D3").Value = "读者的借阅记录分析"
$ws.Range("E3").Value = "读者的系统维护"

# Update the selected cell in the sheet view
$ws.Range("C9").Select()
